$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) "Workblocks" sheet: CloseAllApplications / InitAllApplications
#    no longer contain a workblock by default, at main level
#    (wbCloseAllApplications_*, wbInitAllApplications_*, wbProcess_*)
#    as well as service level (wbCloseAppsRecover_*).
# -----------------------------------------------------------------
$wsWorkblocks = $wb.Worksheets.Item("Workblocks")

# Remove the "CloseApps" (service-level recover) workblock rows (5:6)
$wsWorkblocks.Range("A5:A6").EntireRow.Delete() | Out-Null

# Remove the CloseAllApplications / InitAllApplications / Process
# (main-level) workblock rows - after the deletion above these are
# rows 9:14 (wbCloseAllApplications_Type .. wbProcess_SuppressSuccessful)
$wsWorkblocks.Range("A9:A14").EntireRow.Delete() | Out-Null

# Update the selection shown on this sheet
$wsWorkblocks.Range("A3:C8").Select() | Out-Null

# -----------------------------------------------------------------
# 2) "Tasks" sheet: workflows that contain a workblock no longer pass
#    in_wbType, so the "Task name" columns (SystemTask1_Name /
#    SystemTask2_Name) are removed. The Task1/Task2 dictionary entry
#    no longer needs a description either.
# -----------------------------------------------------------------
$wsTasks = $wb.Worksheets.Item("Tasks")

# Remove "SystemTask1_Name" row (row 3)
$wsTasks.Rows.Item(3).Delete() | Out-Null

# "SystemTask2_Name" is now row 5 - remove it too
$wsTasks.Rows.Item(5).Delete() | Out-Null

# Task2's description cell is no longer populated
$wsTasks.Range("C9").ClearContents() | Out-Null

# Make "Tasks" the active sheet/tab, with the new selection
$wsTasks.Activate() | Out-Null
$wsTasks.Range("B7").Select() | Out-Null

# -----------------------------------------------------------------
# 3) "Constants" sheet keeps its own selection (B5); it simply stops
#    being the active tab now that "Tasks" is selected instead.
# -----------------------------------------------------------------
